$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: Frameworks list - ", Backbone" -> ", Jersey, Angular, Backbone"
# ---------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(", Backbone", $false, $false, $false, $false, $false, $true, 1, $false, ", Jersey, Angular, Backbone", 2)
Write-Host "Edit1 found:" $found

# ---------------------------------------------------------------
# Edit 2: move the "_GoBack" bookmark so it spans from the start of
# the "Languages:" skills bullet through the end of the bullet that
# ends in "PERL". (Adding a bookmark with the same name moves it,
# removing the old occurrence near "Worked" further down.)
# ---------------------------------------------------------------
$rLangStart = $d.Content
$foundLang = $rLangStart.Find.Execute("Languages:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit2 foundLang:" $foundLang

$rPerlEnd = $d.Content
$foundPerl = $rPerlEnd.Find.Execute("C, PERL", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit2 foundPerl:" $foundPerl

$bmRange = $d.Range($rLangStart.Start, $rPerlEnd.End)
Write-Host "Edit2 bmRange text: [" $bmRange.Text "]"
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------
# Edit 3: merge the two runs "           May" + " 2015 " (same
# formatting already) into a single run with the same combined text.
# ---------------------------------------------------------------
$rMay = $d.Content
$foundMay = $rMay.Find.Execute("           May 2015 ", $false, $false, $false, $false, $false, $true, 1, $false, "           May 2015 ", 2)
Write-Host "Edit3 foundMay:" $foundMay

# ---------------------------------------------------------------
# Edit 4: reword the "Worked ..." bullet (Two Sigma) and add a new
# bullet about the recruiting tool underneath it.
# ---------------------------------------------------------------
$rWorked = $d.Content
$foundWorked = $rWorked.Find.Execute(" on the Collaboration, Usability, and Engineering team building workflow optimization tools.", $false, $false, $false, $false, $false, $true, 1, $false, " on the Collaboration, Usability, and Engineering team building the interview scheduling optimizer tool.", 2)
Write-Host "Edit4 foundWorked:" $foundWorked

$rWorkedPara = $d.Content
$foundWorkedPara = $rWorkedPara.Find.Execute("optimizer tool.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit4 foundWorkedPara:" $foundWorkedPara
$workedPara = $rWorkedPara.Paragraphs(1)
$workedParaStart = $workedPara.Range.Start

$workedPara.Range.InsertParagraphAfter()

$allParas = $d.Paragraphs
$newIdx = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -eq $workedParaStart) {
        $newIdx = $i
        break
    }
}
$newPara = $allParas.Item($newIdx + 1)
$newPara.Range.InsertBefore("Developed end to end solution for recruiting by tackling problems with automation, interview distribution, and feedback.")
Write-Host "Edit4 newPara text: [" $newPara.Range.Text "]"

# ---------------------------------------------------------------
# Edit 5: remove the now-obsolete "Developed the frontend website
# and backend database for main product that has around 1500
# users." bullet entirely.
# ---------------------------------------------------------------
$rRemove = $d.Content
$foundRemove = $rRemove.Find.Execute("Developed the frontend website", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Edit5 foundRemove:" $foundRemove
$removePara = $rRemove.Paragraphs(1)
Write-Host "Edit5 removePara text: [" $removePara.Range.Text "]"
$removePara.Range.Delete()
